$wb = $excel.ActiveWorkbook

# --- Fluid_1 sheet ---
$ws1 = $wb.Worksheets.Item("Fluid_1")

# number_compounds: 2 -> 3
$ws1.Range("B3").Value = 3

# fl2 value_1: 0.2 -> 0.1
$ws1.Range("B5").Value = 0.1

# Insert a new row at position 6 (shifts T_in/p_in/m_dot/units down by one)
$ws1.Rows("6:6").Insert()

# New row 6: fl3 (mole fraction entry for Butane)
$ws1.Range("A6").Value = "fl3"
$ws1.Range("B6").Value = 0.1
$ws1.Range("C6").Value = 0.1
$ws1.Range("D6").Value = 0.6
# Copy the shared "True" text (not boolean TRUE) from an existing fixed_1 cell
$ws1.Range("E4").Copy() | Out-Null
$ws1.Range("E6").PasteSpecial() | Out-Null
$ws1.Range("F6").Value = "Butane"
$ws1.Range("G6").Value = "mole fraction"

# Row 7 (was row 6): T_in - value_1 350 -> 370
$ws1.Range("B7").Value = 370

# Row 8 (was row 7): p_in - value_1 1000000 -> 1300000, and now also has fixed_1 = True
$ws1.Range("B8").Value = 1300000
$ws1.Range("E4").Copy() | Out-Null
$ws1.Range("E8").PasteSpecial() | Out-Null

# Row 9 (was row 8): m_dot - value_1 1.2E-2 -> 0.01
$ws1.Range("B9").Value = 0.01

$excel.CutCopyMode = 0

# --- Fluid_2 sheet ---
$ws2 = $wb.Worksheets.Item("Fluid_2")

# m_dot value_1: 3.2000000000000002E-3 -> 9.1999999999999998E-3
$ws2.Range("B7").Value = 0.0092
